# Generate Report for Handoff
# - Overview sheet: mark the e22cf724 file's zh-cn / de-de status as
#   "Ready for handoff" (was "Handed back: in sync with en-US").
# - zh-cn / de-de detail sheets: same Status update for that file's row,
#   plus a refreshed "Latest Handoff Datetime" reflecting the new handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn detail sheet ------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-27 07:55:40"

# --- de-de detail sheet ------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-27 07:55:52"
